$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting (bold, border, alignment) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in team record values for every data row (2-42)
$ws.Range("AD2:AD42").Value = 89
$ws.Range("AE2:AE42").Value = 73
$ws.Range("AF2:AF42").Value = 0
